$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row 3: alternate "xpath" locator for the "header" element.
$ws.Rows("3:3").Insert()
$ws.Range("A3").Value = "header"
$ws.Range("B3").Value = "xpath"
$ws.Range("C3").Value = "xpatch "
$ws.Range("D3").Value = "Serveware & Flatware"

# 2) Insert a new row 8: a highlighted section header before the colour filters.
$ws.Rows("8:8").Insert()
$ws.Range("A8").Value = "collor filters"
$ws.Range("A8").Interior.Color = 65535

# 3) Insert two new rows (10:11): alternate "class name" / "css selector"
#    locators for the same "select Black" element.
$ws.Rows("10:11").Insert()
$ws.Range("A10").Value = "select Black"
$ws.Range("B10").Value = "class name"
$ws.Range("C10").Value = "calss_name"
$ws.Range("D10").Value = "Black"

$ws.Range("A11").Value = "select Black"
$ws.Range("B11").Value = "css selector"
$ws.Range("C11").Value = "css_selector"
$ws.Range("D11").Value = "Black"

# Move the active selection to the first empty row below the table, matching
# where Excel leaves the cursor after these edits.
$ws.Range("A16").Select()
